# moved away from table master_plans, hash codes for folder names
#
# The "plans" worksheet used to source values from a lookup table called
# master_plans; column D ("expand_by_lbd") is no longer produced, so it is
# removed (cells shift left). Row 2's "expand_by" value (old D2 = 40) becomes
# a numeric 40 in the now-merged C2 cell. Every other worksheet simply had
# its previous two-area selection collapsed down to a single active cell.

$wb = $excel.ActiveWorkbook

# --- 1. Plain worksheets: collapse the old "A2:A11 + active cell" selection
#        down to just the active cell. ---
$simpleSelections = @{
    "model_params"      = "A1"
    "dataset_params"    = "E15"
    "transform_factors" = "A1"
    "affine3d"          = "A1"
    "loss_params"       = "A1"
    "plan1"             = "F1"
    "plan2"             = "A1"
    "plan3"             = "A1"
    "plan4"             = "A1"
    "plan9"             = "A1"
    "plan5"             = "A1"
    "plan6"             = "A1"
    "plan7"             = "A1"
    "plan8"             = "A1"
    "plan10"            = "A1"
}

foreach ($name in $simpleSelections.Keys) {
    $sheet = $wb.Worksheets.Item($name)
    $sheet.Activate()
    $sheet.Range($simpleSelections[$name]).Select()
}

# --- 2. "plans" worksheet: drop the old "expand_by_lbd" column (D). ---
$plans = $wb.Worksheets.Item("plans")
$plans.Activate()

# Deleting the whole column shifts every column from E onward one slot to
# the left (E->D, F->E, ... Z->Y).
$plans.Columns("D").Delete()

# Row 2's expand_by (C2) previously held the string "0"; once column D is
# gone the value that used to live in D2 (the number 40) becomes the cell's
# new content, now stored as a genuine number instead of text.
$plans.Range("C2").Value = 40

# Scroll/select so column C is at the left edge and C3 is the active cell.
$plans.Range("C1").Select()
$excel.ActiveWindow.ScrollColumn = 3
$plans.Range("C3").Select()
